$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddAgencyList")

# New column G: header + value, for the Add Agency List sheet
$ws.Range("G1").Value = "AgencyCreatedInAddAgencyFlow"
$ws.Range("G2").Value = "ATMNAgencykpm"

# Column F width adjustment (as captured in the saved workbook)
$ws.Columns.Item(6).ColumnWidth = 13

# Move the active selection to the newly added cell
$ws.Range("G2").Select()
